$wb = $excel.ActiveWorkbook

# Duplicate the R3_P7 sheet (template-style, keeps all number formats/merges/
# column widths/borders identical) and rename/reposition it as R3_P8, right
# after R3_P7 - mirrors how the lab actually built each new box sheet.
$srcSheet = $wb.Worksheets.Item("R3_P7")
$srcSheet.Copy([System.Type]::Missing, $srcSheet)
$ws = $wb.Worksheets.Item($srcSheet.Index + 1)
$ws.Name = "R3_P8"

# Header block
$ws.Range("B1").Value = "2018-07-04"
$ws.Range("B2").Value = "DSPR dead cohort flies"
$ws.Range("B4").Value = "R3_P8"
$ws.Range("B5").Value = "flies that have died after water bath tests"
$ws.Range("B6").Value = "SURF nb#001 pg 101"

# Row 9 - sample labels for cohort 2
$ws.Range("C9").Value = "2p13A4FR"
$ws.Range("D9").Value = "2p5B4HH"
$ws.Range("E9").Value = "2p13B4SHS"
$ws.Range("F9").Value = "2p8A3SHS"
$ws.Range("G9").Value = "2p4A2SHS"
$ws.Range("H9").Value = "2p3A4SHS"
$ws.Range("I9").Value = "2p8A3SR"
$ws.Range("J9").Value = "2p6A1SR"
$ws.Range("K9").ClearContents()

# Row 10 no longer carries the R3_P7 leftover labels
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
